$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-03 Sunday" "2024-11-04 Monday"

Replace-Text "45÷2=22, 1" "69÷8=8, 5"
Replace-Text "63÷8=7, 7" "60÷3=20, 0"
Replace-Text "85÷2=42, 1" "32÷7=4, 4"
Replace-Text "50÷8=6, 2" "54÷7=7, 5"
Replace-Text "15÷6=2, 3" "25÷5=5, 0"

Replace-Text "34÷3=11, 1" "10÷3=3, 1"
Replace-Text "39÷9=4, 3" "26÷7=3, 5"
Replace-Text "45÷5=9, 0" "47÷8=5, 7"
Replace-Text "73÷7=10, 3" "12÷9=1, 3"
Replace-Text "80÷6=13, 2" "10÷7=1, 3"

Replace-Text "71÷2=35, 1" "60÷9=6, 6"
Replace-Text "93÷6=15, 3" "24÷3=8, 0"
Replace-Text "83÷5=16, 3" "93÷5=18, 3"
Replace-Text "31÷7=4, 3" "45÷3=15, 0"
Replace-Text "13÷7=1, 6" "16÷4=4, 0"

Replace-Text "25÷6=4, 1" "92÷9=10, 2"
Replace-Text "95÷2=47, 1" "52÷8=6, 4"
Replace-Text "94÷8=11, 6" "43÷7=6, 1"
Replace-Text "90÷3=30, 0" "47÷4=11, 3"
Replace-Text "50÷7=7, 1" "10÷2=5, 0"

Replace-Text "36÷6=6, 0" "18÷4=4, 2"
Replace-Text "89÷9=9, 8" "86÷4=21, 2"
Replace-Text "24÷7=3, 3" "96÷3=32, 0"
Replace-Text "33÷5=6, 3" "54÷9=6, 0"
Replace-Text "13÷3=4, 1" "89÷3=29, 2"
